$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Abolfazl's row (row 4) degree/grades added across columns F:K
$ws.Range("F4:K4").Value = 100

# Update the active selection to J4 as recorded in the saved view state
$ws.Range("J4").Select()
